$d = $word.ActiveDocument

$replacements = @(
    @{old="64×89="; new="51×36="},
    @{old="13×53="; new="26×80="},
    @{old="31×34="; new="58×43="},
    @{old="59×36="; new="72×58="},
    @{old="25×31="; new="51×38="},
    @{old="99×16="; new="75×85="},
    @{old="22×77="; new="84×51="},
    @{old="44×49="; new="43×25="},
    @{old="92×66="; new="78×40="},
    @{old="67×61="; new="12×90="},
    @{old="34×41="; new="70×84="},
    @{old="17×23="; new="99×33="},
    @{old="60×29="; new="67×83="},
    @{old="21×77="; new="19×64="},
    @{old="70×83="; new="26×47="},
    @{old="62×63="; new="68×85="},
    @{old="82×55="; new="93×60="},
    @{old="33×63="; new="96×62="},
    @{old="36×28="; new="47×14="},
    @{old="27×39="; new="96×86="},
    @{old="73×15="; new="39×36="},
    @{old="75×33="; new="40×81="},
    @{old="50×33="; new="87×76="},
    @{old="52×90="; new="88×73="},
    @{old="15×55="; new="85×42="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
